$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a bare percentage (e.g. "80%") need their
# number format forced to Text first, otherwise Excel auto-converts the typed
# string into a numeric percentage value (0.8 with a Percent format) instead
# of keeping the literal text "80%" the source data uses.
$percentTextCells = @("H3", "H4", "H5", "H7", "H8", "H11", "H12", "H15", "H22", "H29", "H30", "H34")
foreach ($cellRef in $percentTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-05 19:47:53"
$ws.Range("E3").Value = "2026-02-05 19:47:55"
$ws.Range("H3").Value = "80%"
$ws.Range("O3").Value = "-1.8 °C"
$ws.Range("E4").Value = "2026-02-05 19:47:57"
$ws.Range("H4").Value = "67%"
$ws.Range("J4").Value = "989.8 hPa"
$ws.Range("O4").Value = "11.4 °C"
$ws.Range("E5").Value = "2026-02-05 19:48:00"
$ws.Range("H5").Value = "71%"
$ws.Range("J5").Value = "990.0 hPa"
$ws.Range("O5").Value = "9.8 °C"
$ws.Range("E6").Value = "2026-02-05 19:48:02"
$ws.Range("O6").Value = "12.9 °C"
$ws.Range("E7").Value = "2026-02-05 19:48:05"
$ws.Range("H7").Value = "80%"
$ws.Range("J7").Value = "991.6 hPa"
$ws.Range("E8").Value = "2026-02-05 19:48:07"
$ws.Range("H8").Value = "86%"
$ws.Range("O8").Value = "8.7 °C"
$ws.Range("E9").Value = "2026-02-05 19:48:10"
$ws.Range("O9").Value = "2.2 °C"
$ws.Range("E10").Value = "2026-02-05 19:48:12"
$ws.Range("O10").Value = "7.8 °C"
$ws.Range("E11").Value = "2026-02-05 19:48:14"
$ws.Range("H11").Value = "94%"
$ws.Range("J11").Value = "994.8 hPa"
$ws.Range("L11").Value = "40.3 km/h - 273º 19:21 TU"
$ws.Range("M11").Value = "5.1 °C 19:23 TU"
$ws.Range("O11").Value = "0.6 °C"
$ws.Range("E12").Value = "2026-02-05 19:48:17"
$ws.Range("H12").Value = "88%"
$ws.Range("O12").Value = "10.0 °C"
$ws.Range("E13").Value = "2026-02-05 19:48:19"
$ws.Range("O13").Value = "7.8 °C"
$ws.Range("E14").Value = "2026-02-05 19:48:22"
$ws.Range("I14").Value = "7.1 mm"
$ws.Range("L14").Value = "66.6 km/h - 206º 19:23 TU"
$ws.Range("E15").Value = "2026-02-05 19:48:24"
$ws.Range("H15").Value = "82%"
$ws.Range("J15").Value = "990.5 hPa"
$ws.Range("O15").Value = "8.2 °C"
$ws.Range("E16").Value = "2026-02-05 19:48:27"
$ws.Range("O16").Value = "3.8 °C"
$ws.Range("E17").Value = "2026-02-05 19:48:29"
$ws.Range("J17").Value = "995.1 hPa"
$ws.Range("M17").Value = "2.3 °C 19:23 TU"
$ws.Range("E18").Value = "2026-02-05 19:48:32"
$ws.Range("O18").Value = "-4.2 °C"
$ws.Range("E19").Value = "2026-02-05 19:48:34"
$ws.Range("J19").Value = "992.4 hPa"
$ws.Range("E20").Value = "2026-02-05 19:48:37"
$ws.Range("E21").Value = "2026-02-05 19:48:39"
$ws.Range("J21").Value = "990.7 hPa"
$ws.Range("O21").Value = "6.3 °C"
$ws.Range("E22").Value = "2026-02-05 19:48:42"
$ws.Range("H22").Value = "88%"
$ws.Range("O22").Value = "8.7 °C"
$ws.Range("E23").Value = "2026-02-05 19:48:44"
$ws.Range("J23").Value = "989.9 hPa"
$ws.Range("E24").Value = "2026-02-05 19:48:47"
$ws.Range("J24").Value = "989.0 hPa"
$ws.Range("O24").Value = "10.4 °C"
$ws.Range("E25").Value = "2026-02-05 19:48:49"
$ws.Range("J25").Value = "994.0 hPa"
$ws.Range("L25").Value = "20.2 km/h - 208º 19:13 TU"
$ws.Range("M25").Value = "3.9 °C 19:19 TU"
$ws.Range("E26").Value = "2026-02-05 19:48:52"
$ws.Range("E27").Value = "2026-02-05 19:48:54"
$ws.Range("J27").Value = "990.2 hPa"
$ws.Range("O27").Value = "8.7 °C"
$ws.Range("E28").Value = "2026-02-05 19:48:57"
$ws.Range("J28").Value = "993.0 hPa"
$ws.Range("O28").Value = "2.5 °C"
$ws.Range("E29").Value = "2026-02-05 19:48:59"
$ws.Range("H29").Value = "80%"
$ws.Range("O29").Value = "9.1 °C"
$ws.Range("E30").Value = "2026-02-05 19:49:01"
$ws.Range("H30").Value = "66%"
$ws.Range("I30").Value = "5.2 mm"
$ws.Range("E31").Value = "2026-02-05 19:49:04"
$ws.Range("I31").Value = "18.9 mm"
$ws.Range("E32").Value = "2026-02-05 19:49:06"
$ws.Range("O32").Value = "12.0 °C"
$ws.Range("E33").Value = "2026-02-05 19:49:09"
$ws.Range("O33").Value = "9.2 °C"
$ws.Range("E34").Value = "2026-02-05 19:49:11"
$ws.Range("H34").Value = "96%"
$ws.Range("L34").Value = "42.5 km/h - 252º 19:15 TU"
$ws.Range("O34").Value = "3.9 °C"
$ws.Range("E35").Value = "2026-02-05 19:49:14"
$ws.Range("I35").Value = "5.2 mm"
$ws.Range("E36").Value = "2026-02-05 19:49:16"
$ws.Range("K36").Value = "9.6 MJ/m2"
